$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet from "Schedule" to "Sheet1"
$ws.Name = "Sheet1"

# 2. Update the track name in C2
$ws.Range("C2").Value = "1_Brahme Muhurte.mp3"

# 3. Format the header row (A1:D1): thin box border + centered/top-aligned text.
#    Build the combined format on A1 first (so the Borders sub-object batches the
#    border-side writes into a single style change), then propagate that exact
#    style to the remaining header cells via a format-only paste so we don't
#    churn through extra intermediate cell styles.
$headerFirst = $ws.Range("A1")
$headerFirst.HorizontalAlignment = -4108   # xlHAlignCenter
$headerFirst.VerticalAlignment = -4160     # xlVAlignTop
$headerBorders = $headerFirst.Borders
$headerBorders.LineStyle = 1               # xlContinuous
$headerBorders.Weight = 2                  # xlThin

$headerFirst.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)     # xlPasteFormats
$excel.CutCopyMode = $false

Write-Host "Applied schedule.xlsx edits: sheet renamed, C2 updated, header formatted."
